# Insert a new row at row 219 (shifts existing rows 219:289 down to 220:290)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(219).Insert()

# Populate the newly inserted row 219 with the new data record
$ws.Range("A219").Value = 4
$ws.Range("B219").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C219").Value = "Los Lagos"
$ws.Range("D219").Value = 44588
$ws.Range("E219").Value = 10
$ws.Range("F219").Value = 100114013
$ws.Range("G219").Value = "Zanahoria"
$ws.Range("H219").Value = "Sin especificar"
$ws.Range("I219").Value = "Primera"
$ws.Range("J219").Value = 250
$ws.Range("K219").Value = 12500
$ws.Range("L219").Value = 13000
$ws.Range("M219").Value = 12800
$ws.Range("N219").Value = "`$/saco 20 kilos"
$ws.Range("O219").Value = "Región de Ñuble"
$ws.Range("P219").Value = 640
$ws.Range("Q219").Value = 20
$ws.Range("R219").Value = "Hortaliza"
